$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-08 Monday" "2025-09-09 Tuesday"
Replace-Text "571×9=" "548×3="
Replace-Text "216×5=" "439×4="
Replace-Text "863×8=" "987×7="
Replace-Text "431×4=" "483×4="
Replace-Text "838×9=" "486×2="
Replace-Text "260×3=" "878×8="
Replace-Text "822×7=" "983×2="
Replace-Text "523×3=" "281×5="
Replace-Text "180×8=" "671×3="
Replace-Text "827×8=" "282×2="
Replace-Text "123×8=" "852×4="
Replace-Text "414×5=" "952×2="
Replace-Text "601×5=" "597×4="
Replace-Text "642×6=" "977×5="
Replace-Text "874×5=" "977×3="
Replace-Text "282×5=" "280×5="
Replace-Text "349×9=" "709×6="
Replace-Text "975×6=" "549×8="
Replace-Text "224×5=" "932×9="
Replace-Text "213×2=" "685×5="
Replace-Text "834×3=" "320×9="
Replace-Text "877×7=" "907×7="
Replace-Text "496×4=" "637×2="
Replace-Text "721×6=" "477×3="
Replace-Text "891×9=" "995×3="
